$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9609749913215637
$ws.Range("B1").Value = 2.179699659347534
$ws.Range("C1").Value = 8.268882751464844
$ws.Range("D1").Value = 1.888365030288696
$ws.Range("E1").Value = 1.175182461738586
